$d = $word.ActiveDocument

# 1) Remove bold from the "Cork" run (keeps bCs) and
# 2) split the following "  is the second largest city in " run into
#    "  i" / "s the second largest city in " (artifact of the retype),
#    without changing any other run properties.
$rCork = $d.Range(0, 4)
$rCork.Font.Bold = 0

$rSplit = $d.Range(4, 7)
$rSplit.Font.Bold = 1
$rSplit.Font.Bold = 0

# 3) Move the _GoBack bookmark into its own paragraph at the end of the
#    document, and drop the now-redundant trailing empty paragraph.
$pLast = $d.Paragraphs(13)
$delRange = $d.Range($pLast.Range.End - 1, $pLast.Range.End)
$delRange.Delete()

$pMerged = $d.Paragraphs(13)
$endPos = $pMerged.Range.End
$pilcrow = $d.Range($endPos - 1, $endPos)
$pilcrow.InsertBefore("`r")
